$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 10 ("case_when(...)") entirely; this shifts rows 11-17 up to 10-16
$ws.Rows.Item(10).Delete()
